$d = $word.ActiveDocument

$d.Content.Find.Execute("57×54=3078", $true, $false, $false, $false, $false, $true, 1, $false, "30×37=1110", 2)
$d.Content.Find.Execute("28×33=924", $true, $false, $false, $false, $false, $true, 1, $false, "13×83=1079", 2)
$d.Content.Find.Execute("46×84=3864", $true, $false, $false, $false, $false, $true, 1, $false, "96×23=2208", 2)
$d.Content.Find.Execute("42×16=672", $true, $false, $false, $false, $false, $true, 1, $false, "29×68=1972", 2)
$d.Content.Find.Execute("17×39=663", $true, $false, $false, $false, $false, $true, 1, $false, "61×53=3233", 2)
$d.Content.Find.Execute("72×20=1440", $true, $false, $false, $false, $false, $true, 1, $false, "90×67=6030", 2)
$d.Content.Find.Execute("51×43=2193", $true, $false, $false, $false, $false, $true, 1, $false, "45×35=1575", 2)
$d.Content.Find.Execute("19×19=361", $true, $false, $false, $false, $false, $true, 1, $false, "85×18=1530", 2)
$d.Content.Find.Execute("53×31=1643", $true, $false, $false, $false, $false, $true, 1, $false, "48×32=1536", 2)
$d.Content.Find.Execute("56×98=5488", $true, $false, $false, $false, $false, $true, 1, $false, "86×71=6106", 2)
$d.Content.Find.Execute("89×60=5340", $true, $false, $false, $false, $false, $true, 1, $false, "86×41=3526", 2)
$d.Content.Find.Execute("31×18=558", $true, $false, $false, $false, $false, $true, 1, $false, "59×85=5015", 2)
$d.Content.Find.Execute("67×21=1407", $true, $false, $false, $false, $false, $true, 1, $false, "54×76=4104", 2)
$d.Content.Find.Execute("63×89=5607", $true, $false, $false, $false, $false, $true, 1, $false, "50×87=4350", 2)
$d.Content.Find.Execute("52×41=2132", $true, $false, $false, $false, $false, $true, 1, $false, "17×95=1615", 2)
$d.Content.Find.Execute("48×90=4320", $true, $false, $false, $false, $false, $true, 1, $false, "34×72=2448", 2)
$d.Content.Find.Execute("40×37=1480", $true, $false, $false, $false, $false, $true, 1, $false, "45×13=585", 2)
$d.Content.Find.Execute("91×56=5096", $true, $false, $false, $false, $false, $true, 1, $false, "39×42=1638", 2)
$d.Content.Find.Execute("81×22=1782", $true, $false, $false, $false, $false, $true, 1, $false, "78×45=3510", 2)
$d.Content.Find.Execute("49×65=3185", $true, $false, $false, $false, $false, $true, 1, $false, "33×98=3234", 2)
$d.Content.Find.Execute("96×26=2496", $true, $false, $false, $false, $false, $true, 1, $false, "26×56=1456", 2)
$d.Content.Find.Execute("99×39=3861", $true, $false, $false, $false, $false, $true, 1, $false, "14×15=210", 2)
$d.Content.Find.Execute("30×88=2640", $true, $false, $false, $false, $false, $true, 1, $false, "71×25=1775", 2)
$d.Content.Find.Execute("94×22=2068", $true, $false, $false, $false, $false, $true, 1, $false, "69×21=1449", 2)
$d.Content.Find.Execute("48×61=2928", $true, $false, $false, $false, $false, $true, 1, $false, "16×84=1344", 2)
